$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set all Fitness values (column C, rows 2-252) to 7586
$ws.Range("C2:C252").Value = 7586
